$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 8 - Nuri (Muhasebe) : add "Çalıştığı Yıl" / "İkramiye" / "Durum1" / "Durum2"
# ---------------------------------------------------------------------------
$ws.Range("E8").Formula = "=YEAR(I1)-YEAR(D8)"
$ws.Range("G8").Formula = "=PRODUCT(F8*0.25)"
$ws.Range("H8").Value = 10
$ws.Range("I8").NumberFormat = "#,##0.00"
$ws.Range("I8").Value = 1111.5

# ---------------------------------------------------------------------------
# Row 9 - Kerim (Bilgi İşlem)
# ---------------------------------------------------------------------------
$ws.Range("E9").Formula = "=YEAR(I1)-YEAR(D9)"
$ws.Range("G9").Formula = "=PRODUCT(F9*0.15)"
$ws.Range("H9").Value = 7
$ws.Range("I9").Value = 877.5

# ---------------------------------------------------------------------------
# Row 10 - Ergün (Bilgi İşlem)
# ---------------------------------------------------------------------------
$ws.Range("E10").Value = 3
$ws.Range("G10").Formula = "=PRODUCT(950*0.25)"
$ws.Range("H10").Value = 7
$ws.Range("I10").NumberFormat = "#,##0.00"
$ws.Range("I10").Value = 1111.5

# ---------------------------------------------------------------------------
# Row 11 - Erhan (Bilgi İşlem)
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 4
$ws.Range("G11").Formula = "=PRODUCT(1150*0.25)"
$ws.Range("H11").Value = 7
$ws.Range("I11").NumberFormat = "#,##0.00"
$ws.Range("I11").Value = 1299.5

# ---------------------------------------------------------------------------
# Row 12 - Celal (İnsan Kaynakları)
# ---------------------------------------------------------------------------
$ws.Range("E12").Value = 6
$ws.Range("G12").Formula = "=PRODUCT(F12*0.4)"
$ws.Range("H12").Value = 15
$ws.Range("I12").Value = 1536

# ---------------------------------------------------------------------------
# Row 13 separator - slightly taller
# ---------------------------------------------------------------------------
$ws.Rows("13").RowHeight = 3.1

# ---------------------------------------------------------------------------
# Personnel-form block (Numara / Ad Soyad / Bölüm)
# ---------------------------------------------------------------------------
$ws.Range("L14").Value = 20215070055
$ws.Range("L15").Value = "Muhammed Ali Harmancı"
$ws.Range("L16").Value = "Yönetim Bilişim Sistemleri"

# ---------------------------------------------------------------------------
# View state - select the department cell (last thing the author touched)
# ---------------------------------------------------------------------------
$ws.Range("L16:N16").Select()
